# Updated simulation sweep results (simplify data added to simOut)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column widths shrank slightly to fit the refreshed (shorter) numeric values
$ws.Columns.Item(1).ColumnWidth = 21.5
$ws.Columns.Item(2).ColumnWidth = 29.5
$ws.Columns.Item(3).ColumnWidth = 28
$ws.Columns.Item(4).ColumnWidth = 34.333333333333336
$ws.Columns.Item(5).ColumnWidth = 28.166666666666668
$ws.Columns.Item(6).ColumnWidth = 26.833333333333332
$ws.Columns.Item(7).ColumnWidth = 33.166666666666664
$ws.Columns.Item(8).ColumnWidth = 28.333333333333332
$ws.Columns.Item(9).ColumnWidth = 27.166666666666668

$ws.Range("B2").Value = 34.749984565889072
$ws.Range("C2").Value = 16.027423637899854
$ws.Range("D2").Value = 0.46122102896219791
$ws.Range("E2").Value = 30.65946564110218
$ws.Range("F2").Value = 15.315465106263842
$ws.Range("G2").Value = 0.49953463917296326
$ws.Range("H2").Value = 337.5
$ws.Range("I2").Value = 301

$ws.Range("B3").Value = 34.792275828140532
$ws.Range("C3").Value = 16.070862018274894
$ws.Range("D3").Value = 0.46190890465626089
$ws.Range("E3").Value = 30.69134123498948
$ws.Range("F3").Value = 15.331044661654182
$ws.Range("G3").Value = 0.49952345009204474
$ws.Range("H3").Value = 337.5
$ws.Range("I3").Value = 301

$ws.Range("B4").Value = 34.822369520987223
$ws.Range("C4").Value = 16.103495876470845
$ws.Range("D4").Value = 0.46244687245551652
$ws.Range("E4").Value = 30.736419501104685
$ws.Range("F4").Value = 15.386645847454201
$ws.Range("G4").Value = 0.50059981276938248
$ws.Range("H4").Value = 338
$ws.Range("I4").Value = 301

$ws.Range("B5").Value = 34.84296501344155
$ws.Range("C5").Value = 16.121617767400551
$ws.Range("D5").Value = 0.46269362441403111
$ws.Range("E5").Value = 30.759968894421611
$ws.Range("F5").Value = 15.387899196629453
$ws.Range("G5").Value = 0.50025730680826808
$ws.Range("H5").Value = 338
$ws.Range("I5").Value = 301

$ws.Range("B6").Value = 34.862395989746702
$ws.Range("C6").Value = 16.128958479658881
$ws.Range("D6").Value = 0.46264629902094312
$ws.Range("E6").Value = 30.820760503377322
$ws.Range("F6").Value = 15.457196285304263
$ws.Range("G6").Value = 0.50151897723647909
$ws.Range("H6").Value = 338
$ws.Range("I6").Value = 300.5

$ws.Range("B7").Value = 34.868843693494135
$ws.Range("C7").Value = 16.137677079980662
$ws.Range("D7").Value = 0.46281078953563481
$ws.Range("E7").Value = 30.851928732119447
$ws.Range("F7").Value = 15.487551939028245
$ws.Range("G7").Value = 0.50199623088408096
$ws.Range("H7").Value = 338
$ws.Range("I7").Value = 300.5

$ws.Range("B8").Value = 34.877559193744702
$ws.Range("C8").Value = 16.145021761430456
$ws.Range("D8").Value = 0.46290572318277562
$ws.Range("E8").Value = 30.870098301379386
$ws.Range("F8").Value = 15.505557011451746
$ws.Range("G8").Value = 0.50228401801878619
$ws.Range("H8").Value = 338
$ws.Range("I8").Value = 300.5

$ws.Range("B9").Value = 34.890879815397469
$ws.Range("C9").Value = 16.13160958508314
$ws.Range("D9").Value = 0.46234459177966053
$ws.Range("E9").Value = 30.890394099833877
$ws.Range("F9").Value = 15.52067963950469
$ws.Range("G9").Value = 0.50244356188347106
$ws.Range("H9").Value = 338
$ws.Range("I9").Value = 300.5

$ws.Range("B10").Value = 34.903941108611413
$ws.Range("C10").Value = 16.11908090092891
$ws.Range("D10").Value = 0.46181263172461778
$ws.Range("E10").Value = 30.904004155888181
$ws.Range("F10").Value = 15.543904029220483
$ws.Range("G10").Value = 0.50297378782415425
$ws.Range("H10").Value = 338
$ws.Range("I10").Value = 300.5

$ws.Range("B11").Value = 34.918023851338937
$ws.Range("C11").Value = 16.133717676544649
$ws.Range("D11").Value = 0.46204555404489189
$ws.Range("E11").Value = 30.924690051711572
$ws.Range("F11").Value = 15.537925766122775
$ws.Range("G11").Value = 0.50244402579752956
$ws.Range("H11").Value = 338
$ws.Range("I11").Value = 300.5

Write-Host "Values updated"
